$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL (B2): bump ig path to v02
$wsMeta.Range("B2").Value = "https://ncpi-fhir.github.io/ncpi-fhir-ig-v02/ValueSet/condition-inheritance-vs"

# Version (B3): 0.1.0 -> 0.2.0
$wsMeta.Range("B3").Value = "0.2.0"

# Date (B8): update publish date/time
$wsMeta.Range("B8").Value = "2022-05-26T12:34:56-05:00"

# --- "Include from Condition Inheri" sheet updates ---
$wsInclude = $wb.Worksheets.Item("Include from Condition Inheri")

# System URI (B4): bump ig path to v02
$wsInclude.Range("B4").Value = "https://ncpi-fhir.github.io/ncpi-fhir-ig-v02/CodeSystem/ConditionInheritanceMode"
